# Update column G ("K") values for rows 2-12 on Sheet1.
# Regenerated save_data: K (strikeouts) replaces old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 3
    3  = 3
    4  = 2
    5  = 2
    6  = 3
    7  = 3
    8  = 5
    9  = 8
    10 = 6
    11 = 5
    12 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
